$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 32.77329241352902
$ws.Range("C2").Value = 16.274835892552474
$ws.Range("D2").Value = 0.49658837101865666
$ws.Range("E2").Value = 29.320677753896625
$ws.Range("F2").Value = 14.579297759489732
$ws.Range("G2").Value = 0.49723604214954376
$ws.Range("H2").Value = 299.05770891393894
$ws.Range("I2").Value = 265.42554296131311
$ws.Range("B3").Value = 32.774252026261962
$ws.Range("C3").Value = 16.275195457947795
$ws.Range("D3").Value = 0.49658480214610251
$ws.Range("E3").Value = 29.272639081544099
$ws.Range("F3").Value = 14.483270274595737
$ws.Range("G3").Value = 0.49477159316760039
$ws.Range("H3").Value = 299
$ws.Range("I3").Value = 265.50939299472316
$ws.Range("B4").Value = 32.772793603631591
$ws.Range("C4").Value = 16.275607717764231
$ws.Range("D4").Value = 0.49661948000553457
$ws.Range("E4").Value = 29.316903451891847
$ws.Range("F4").Value = 14.524171121694941
$ws.Range("G4").Value = 0.49541968665035402
$ws.Range("H4").Value = 299
$ws.Range("I4").Value = 265.68202053819147
$ws.Range("B5").Value = 32.772481962142614
$ws.Range("C5").Value = 16.276321602329066
$ws.Range("D5").Value = 0.49664598552928596
$ws.Range("E5").Value = 29.30691429262793
$ws.Range("F5").Value = 14.603291602409788
$ws.Range("G5").Value = 0.49828826933455783
$ws.Range("H5").Value = 298.9497081585269
$ws.Range("I5").Value = 265.39255207683027
$ws.Range("B6").Value = 32.769371324296081
$ws.Range("C6").Value = 16.280739717284661
$ws.Range("D6").Value = 0.49682795425537163
$ws.Range("E6").Value = 29.297991960860699
$ws.Range("F6").Value = 14.547901135108079
$ws.Range("G6").Value = 0.49654942750147096
$ws.Range("H6").Value = 299.07415270253398
$ws.Range("I6").Value = 265.36600548902254
$ws.Range("B7").Value = 32.76122057453825
$ws.Range("C7").Value = 16.285613400966337
$ws.Range("D7").Value = 0.49710032518212649
$ws.Range("E7").Value = 29.319446752684335
$ws.Range("F7").Value = 14.528667127739583
$ws.Range("G7").Value = 0.49553005724466526
$ws.Range("H7").Value = 298.94012124375212
$ws.Range("I7").Value = 265.59005044578947
$ws.Range("B8").Value = 32.794728442994554
$ws.Range("C8").Value = 16.298977203582815
$ws.Range("D8").Value = 0.4969999136268049
$ws.Range("E8").Value = 29.400424879797164
$ws.Range("F8").Value = 14.528099914219124
$ws.Range("G8").Value = 0.4941459170612964
$ws.Range("H8").Value = 299.04531566934696
$ws.Range("I8").Value = 266.0247110216535
$ws.Range("B9").Value = 32.811284433865559
$ws.Range("C9").Value = 16.342171135941747
$ws.Range("D9").Value = 0.49806557158349091
$ws.Range("E9").Value = 29.509394935202451
$ws.Range("F9").Value = 14.554407235259344
$ws.Range("G9").Value = 0.49321266217820853
$ws.Range("H9").Value = 299.5
$ws.Range("I9").Value = 266.9838460207726
$ws.Range("B10").Value = 32.883544459803936
$ws.Range("C10").Value = 16.376378724027308
$ws.Range("D10").Value = 0.49801136079005731
$ws.Range("E10").Value = 29.773284346783221
$ws.Range("F10").Value = 14.574585052800272
$ws.Range("G10").Value = 0.48951888824367962
$ws.Range("H10").Value = 299.52225066632963
$ws.Range("I10").Value = 268.96354852795059
$ws.Range("B11").Value = 33.010550936512153
$ws.Range("C11").Value = 16.238632654447056
$ws.Range("D11").Value = 0.49192249731542365
$ws.Range("E11").Value = 30.095801120460571
$ws.Range("F11").Value = 14.949946014553548
$ws.Range("G11").Value = 0.49674524212581461
$ws.Range("H11").Value = 300.92703705314574
$ws.Range("I11").Value = 271.99088509064114
